$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: Unmerge the big header merge-ranges that need to be split so that
# the previously-empty interior cells (AF/AZ columns) become independently
# addressable anchor cells.
# ---------------------------------------------------------------------------
$ws.Range("W1:AO1").UnMerge()
$ws.Range("AP1:BH1").UnMerge()
$ws.Range("W2:AO2").UnMerge()
$ws.Range("AP2:BH2").UnMerge()
$ws.Range("W3:AO3").UnMerge()
$ws.Range("AP3:BH3").UnMerge()
$ws.Range("W5:AO5").UnMerge()
$ws.Range("AP5:BH5").UnMerge()

# ---------------------------------------------------------------------------
# Step 2: Populate the newly freed cells. Each of these header cells needs
# the same "boxed, centered, bold Times-New-Roman" look already used by the
# header row (row 6), plus a colored fill matching the corresponding section
# (green for FT-x rows, orange for CO-x rows, plain for THEORY, grey for the
# "Question numbers mapping" rows). Copying borders/font/alignment from an
# *unmerged* cell (A6) avoids the border artefacts that occur when copying
# format directly from a cell that is part of a merged range.
# ---------------------------------------------------------------------------
$ws.Range("A6").Copy()
$ws.Range("AF1").PasteSpecial(-4122)
$ws.Range("AZ1").PasteSpecial(-4122)
$ws.Range("AF2").PasteSpecial(-4122)
$ws.Range("AZ2").PasteSpecial(-4122)
$ws.Range("AF3").PasteSpecial(-4122)
$ws.Range("AZ3").PasteSpecial(-4122)
$ws.Range("AF5").PasteSpecial(-4122)
$ws.Range("AZ5").PasteSpecial(-4122)

# Row 1 fill = green (matches FT-I / FT-III header cells)
$ws.Range("AF1").Interior.Color = 0x50D092
$ws.Range("AZ1").Interior.Color = 0x50D092
$ws.Range("AF1").Value = "FT-II , FT-IV"
$ws.Range("AZ1").Value = "FT-III, REPORT"

# Row 2 fill = orange (matches CO1 / CO2 header cells)
$ws.Range("AF2").Interior.Color = 0x00C3FF
$ws.Range("AZ2").Interior.Color = 0x00C3FF
$ws.Range("AF2").Value = "CO3"
$ws.Range("AP2").Value = "CO4"
$ws.Range("AZ2").Value = "CO5"

# Row 3 (no fill, matches the other THEORY cells); also normalize the THEORY
# text (the multi-line text becomes a single line without embedded newlines)
$theoryText = "THEORY (for either/or Q, award marks for the attempted students only)"
$ws.Range("D3").Value = $theoryText
$ws.Range("W3").Value = $theoryText
$ws.Range("AF3").Value = $theoryText
$ws.Range("AP3").Value = $theoryText
$ws.Range("AZ3").Value = $theoryText

# Row 4: numeric MAX MARKS values get reshuffled between columns
$ws.Range("AC4").Value = 8
$ws.Range("AD4").Value = 8
$ws.Range("AE4").Value = 15
$ws.Range("AH4").Value = 1
$ws.Range("AI4").Value = 1
$ws.Range("AJ4").Value = 1
$ws.Range("AN4").Value = 8
$ws.Range("AV4").Value = 8
$ws.Range("AW4").Value = 8
$ws.Range("AX4").Value = 8
$ws.Range("AY4").Value = 8
$ws.Range("BA4").Value = 1
$ws.Range("BB4").Value = 1
$ws.Range("BC4").Value = 1
$ws.Range("BD4").Value = 1

# Row 5 fill = grey (matches the other "Question numbers mapping" cells)
$ws.Range("AF5").Interior.Color = 0xC0C0C0
$ws.Range("AZ5").Interior.Color = 0xC0C0C0
$ws.Range("AF5").Value = "Question numbers mapping"
$ws.Range("AZ5").Value = "Question numbers mapping"

# Row 6: re-map the question-number header labels for columns AC..BD
$ws.Range("AC6").Value = "Q12.A"
$ws.Range("AD6").Value = "Q12.B"
$ws.Range("AE6").Value = "Q15"
$ws.Range("AF6").Value = "Q7"
$ws.Range("AG6").Value = "Q8"
$ws.Range("AH6").Value = "Q9"
$ws.Range("AI6").Value = "Q10"
$ws.Range("AJ6").Value = "Q11"
$ws.Range("AK6").Value = "Q13.A"
$ws.Range("AL6").Value = "Q13.B"
$ws.Range("AM6").Value = "Q14.A"
$ws.Range("AN6").Value = "Q14.B"
$ws.Range("AO6").Value = "Q16"
$ws.Range("AV6").Value = "Q12.A"
$ws.Range("AW6").Value = "Q12.B"
$ws.Range("AX6").Value = "Q13.A"
$ws.Range("AY6").Value = "Q13.B"
$ws.Range("AZ6").Value = "Q7"
$ws.Range("BA6").Value = "Q8"
$ws.Range("BB6").Value = "Q9"
$ws.Range("BC6").Value = "Q10"
$ws.Range("BD6").Value = "Q11"

# ---------------------------------------------------------------------------
# Step 3: Re-merge the header ranges using the new, split layout.
# ---------------------------------------------------------------------------
$ws.Range("W1:AE1").Merge()
$ws.Range("AF1:AO1").Merge()
$ws.Range("AP1:AY1").Merge()
$ws.Range("AZ1:BH1").Merge()

$ws.Range("W2:AE2").Merge()
$ws.Range("AF2:AO2").Merge()
$ws.Range("AP2:AY2").Merge()
$ws.Range("AZ2:BH2").Merge()

$ws.Range("W3:AE3").Merge()
$ws.Range("AF3:AO3").Merge()
$ws.Range("AP3:AY3").Merge()
$ws.Range("AZ3:BH3").Merge()

$ws.Range("W5:AE5").Merge()
$ws.Range("AF5:AO5").Merge()
$ws.Range("AP5:AY5").Merge()
$ws.Range("AZ5:BH5").Merge()

Write-Output "edit complete"
